# Adding basics of Gully Sediment Contribution
#
# Insert a new "Gullies" worksheet between the existing "PlanningUnits" and
# "testAscTable" sheets, and populate it with gully erosion data keyed by
# planning unit.

$wb = $excel.ActiveWorkbook

# "PlanningUnits" is the first sheet; inserting After it places the new
# sheet immediately before "testAscTable", matching the target layout:
#   PlanningUnits, Gullies, testAscTable
$planningUnits = $wb.Worksheets.Item(1)
$gullies = $wb.Worksheets.Add($null, $planningUnits)
$gullies.Name = "Gullies"

# Header row. Write column C before column B so new shared-string entries
# are appended in the order GullyErosionVolume, then PlanningUnit.
$gullies.Range("A1").Value = "Identifier"
$gullies.Range("C1").Value = "GullyErosionVolume"
$gullies.Range("B1").Value = "PlanningUnit"
$gullies.Range("D1").Value = "ChannelLength"

# Identifier, PlanningUnit, GullyErosionVolume, ChannelLength(/Area column)
$data = @(
    @(1, 3, 341207.89779999998, 5730),
    @(2, 3, 12978.785879999999, 502),
    @(3, 1, 46200, 643),
    @(4, 2, 40639.491300000002, 2044),
    @(5, 5, 40051, 1961)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $gullies.Cells.Item($row, 1).Value = $data[$i][0]
    $gullies.Cells.Item($row, 2).Value = $data[$i][1]
    $gullies.Cells.Item($row, 3).Value = $data[$i][2]
    $gullies.Cells.Item($row, 4).Value = $data[$i][3]
}
